# Populate the new "Speculative % Cover" sheet with data, mirroring the
# Biomass sheet row formatting (block color-coding) and labels, excluding
# "Weeds" entries, per the commit: speculative coverage cleaned and all
# control treatments combined.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Biomass")
$ws2 = $wb.Worksheets.Item("Speculative % Cover")

# Narrow column B on Biomass slightly (cosmetic width tweak from the edit)
$ws1.Columns.Item(2).ColumnWidth = 16

# Header row
$ws1.Range("A1").Copy() | Out-Null
$ws2.Range("A1").PasteSpecial(-4122) | Out-Null
$ws2.Range("A1").Value = "Barcode"
$ws2.Range("B1").Value = "Speculative % Cover"

$ws1.Range("A2").Copy() | Out-Null
$ws2.Range("A2").PasteSpecial(-4122) | Out-Null
$ws2.Range("A2").Value = "FSB-1_N_Control"
$ws2.Range("B2").Value = 0

$ws1.Range("A2").Copy() | Out-Null
$ws2.Range("A3").PasteSpecial(-4122) | Out-Null
$ws2.Range("A3").Value = "FSB-1_C_Rapeseed"
$ws2.Range("B3").Value = 0

$ws1.Range("A2").Copy() | Out-Null
$ws2.Range("A4").PasteSpecial(-4122) | Out-Null
$ws2.Range("A4").Value = "FSB-1_H_Oriental Mustard"
$ws2.Range("B4").Value = 40

$ws1.Range("A2").Copy() | Out-Null
$ws2.Range("A5").PasteSpecial(-4122) | Out-Null
$ws2.Range("A5").Value = "FSB-1_F_Mustard"
$ws2.Range("B5").Value = 45

$ws1.Range("A2").Copy() | Out-Null
$ws2.Range("A6").PasteSpecial(-4122) | Out-Null
$ws2.Range("A6").Value = "FSB-1_G_Brown Mustard"
$ws2.Range("B6").Value = 50

$ws1.Range("A2").Copy() | Out-Null
$ws2.Range("A7").PasteSpecial(-4122) | Out-Null
$ws2.Range("A7").Value = "FSB-1_A_Radish"
$ws2.Range("B7").Value = 10

$ws1.Range("A2").Copy() | Out-Null
$ws2.Range("A8").PasteSpecial(-4122) | Out-Null
$ws2.Range("A8").Value = "FSB-1_I_Arugula"
$ws2.Range("B8").Value = 2.5

$ws1.Range("A2").Copy() | Out-Null
$ws2.Range("A9").PasteSpecial(-4122) | Out-Null
$ws2.Range("A9").Value = "FSB-1_L_Red Clover"
$ws2.Range("B9").Value = 2.5

$ws1.Range("A2").Copy() | Out-Null
$ws2.Range("A10").PasteSpecial(-4122) | Out-Null
$ws2.Range("A10").Value = "FSB-1_D_Collard"
$ws2.Range("B10").Value = 45

$ws1.Range("A2").Copy() | Out-Null
$ws2.Range("A11").PasteSpecial(-4122) | Out-Null
$ws2.Range("A11").Value = "FSB-1_B_Turnip"
$ws2.Range("B11").Value = 35

$ws1.Range("A2").Copy() | Out-Null
$ws2.Range("A12").PasteSpecial(-4122) | Out-Null
$ws2.Range("A12").Value = "FSB-1_K_Kale"
$ws2.Range("B12").Value = 15

$ws1.Range("A2").Copy() | Out-Null
$ws2.Range("A13").PasteSpecial(-4122) | Out-Null
$ws2.Range("A13").Value = "FSB-1_M_Alfalfa"
$ws2.Range("B13").Value = 10

$ws1.Range("A2").Copy() | Out-Null
$ws2.Range("A14").PasteSpecial(-4122) | Out-Null
$ws2.Range("A14").Value = "FSB-1_O_Control"
$ws2.Range("B14").Value = 0

$ws1.Range("A2").Copy() | Out-Null
$ws2.Range("A15").PasteSpecial(-4122) | Out-Null
$ws2.Range("A15").Value = "FSB-1_J_Winter Camelina"
$ws2.Range("B15").Value = 0

$ws1.Range("A2").Copy() | Out-Null
$ws2.Range("A16").PasteSpecial(-4122) | Out-Null
$ws2.Range("A16").Value = "FSB-1_P_Control"
$ws2.Range("B16").Value = 0

$ws1.Range("A2").Copy() | Out-Null
$ws2.Range("A17").PasteSpecial(-4122) | Out-Null
$ws2.Range("A17").Value = "FSB-1_E_Broadleaf Mustard"
$ws2.Range("B17").Value = 35

$ws1.Range("A31").Copy() | Out-Null
$ws2.Range("A18").PasteSpecial(-4122) | Out-Null
$ws2.Range("A18").Value = "FSB-2_D_Collard"
$ws2.Range("B18").Value = 55

$ws1.Range("A31").Copy() | Out-Null
$ws2.Range("A19").PasteSpecial(-4122) | Out-Null
$ws2.Range("A19").Value = "FSB-2_B_Turnip"
$ws2.Range("B19").Value = 20

$ws1.Range("A31").Copy() | Out-Null
$ws2.Range("A20").PasteSpecial(-4122) | Out-Null
$ws2.Range("A20").Value = "FSB-2_K_Kale"
$ws2.Range("B20").Value = 15

$ws1.Range("A31").Copy() | Out-Null
$ws2.Range("A21").PasteSpecial(-4122) | Out-Null
$ws2.Range("A21").Value = "FSB-2_L_Red Clover"
$ws2.Range("B21").Value = 2.5

$ws1.Range("A31").Copy() | Out-Null
$ws2.Range("A22").PasteSpecial(-4122) | Out-Null
$ws2.Range("A22").Value = "FSB-3_P_Control"
$ws2.Range("B22").Value = 0

$ws1.Range("A31").Copy() | Out-Null
$ws2.Range("A23").PasteSpecial(-4122) | Out-Null
$ws2.Range("A23").Value = "FSB-2_J_Winter Camelina"
$ws2.Range("B23").Value = 0

$ws1.Range("A31").Copy() | Out-Null
$ws2.Range("A24").PasteSpecial(-4122) | Out-Null
$ws2.Range("A24").Value = "FSB-2_H_Oriental Mustard"
$ws2.Range("B24").Value = 45

$ws1.Range("A31").Copy() | Out-Null
$ws2.Range("A25").PasteSpecial(-4122) | Out-Null
$ws2.Range("A25").Value = "FSB-2_F_Mustard"
$ws2.Range("B25").Value = 40

$ws1.Range("A31").Copy() | Out-Null
$ws2.Range("A26").PasteSpecial(-4122) | Out-Null
$ws2.Range("A26").Value = "FSB-2_O_Control"
$ws2.Range("B26").Value = 0

$ws1.Range("A31").Copy() | Out-Null
$ws2.Range("A27").PasteSpecial(-4122) | Out-Null
$ws2.Range("A27").Value = "FSB-2_C_Rapeseed"
$ws2.Range("B27").Value = 7.5

$ws1.Range("A31").Copy() | Out-Null
$ws2.Range("A28").PasteSpecial(-4122) | Out-Null
$ws2.Range("A28").Value = "FSB-2_I_Arugula"
$ws2.Range("B28").Value = 5

$ws1.Range("A31").Copy() | Out-Null
$ws2.Range("A29").PasteSpecial(-4122) | Out-Null
$ws2.Range("A29").Value = "FSB-2_G_Brown Mustard"
$ws2.Range("B29").Value = 35

$ws1.Range("A31").Copy() | Out-Null
$ws2.Range("A30").PasteSpecial(-4122) | Out-Null
$ws2.Range("A30").Value = "FSB-2_M_Alfalfa"
$ws2.Range("B30").Value = 15

$ws1.Range("A31").Copy() | Out-Null
$ws2.Range("A31").PasteSpecial(-4122) | Out-Null
$ws2.Range("A31").Value = "FSB-2_N_Control"
$ws2.Range("B31").Value = 0

$ws1.Range("A31").Copy() | Out-Null
$ws2.Range("A32").PasteSpecial(-4122) | Out-Null
$ws2.Range("A32").Value = "FSB-2_A_Radish"
$ws2.Range("B32").Value = 2.5

$ws1.Range("A31").Copy() | Out-Null
$ws2.Range("A33").PasteSpecial(-4122) | Out-Null
$ws2.Range("A33").Value = "FSB-2_E_Broadleaf Mustard"
$ws2.Range("B33").Value = 25

$ws1.Range("A60").Copy() | Out-Null
$ws2.Range("A34").PasteSpecial(-4122) | Out-Null
$ws2.Range("A34").Value = "FSB-3_A_Radish"
$ws2.Range("B34").Value = 5

$ws1.Range("A60").Copy() | Out-Null
$ws2.Range("A35").PasteSpecial(-4122) | Out-Null
$ws2.Range("A35").Value = "FSB-3_C_Rapeseed"
$ws2.Range("B35").Value = 0

$ws1.Range("A60").Copy() | Out-Null
$ws2.Range("A36").PasteSpecial(-4122) | Out-Null
$ws2.Range("A36").Value = "FSB-3_N_Control"
$ws2.Range("B36").Value = 0

$ws1.Range("A60").Copy() | Out-Null
$ws2.Range("A37").PasteSpecial(-4122) | Out-Null
$ws2.Range("A37").Value = "FSB-3_J_Winter Camelina"
$ws2.Range("B37").Value = 0

$ws1.Range("A60").Copy() | Out-Null
$ws2.Range("A38").PasteSpecial(-4122) | Out-Null
$ws2.Range("A38").Value = "FSB-3_E_Broadleaf Mustard"
$ws2.Range("B38").Value = 15

$ws1.Range("A60").Copy() | Out-Null
$ws2.Range("A39").PasteSpecial(-4122) | Out-Null
$ws2.Range("A39").Value = "FSB-3_B_Turnip"
$ws2.Range("B39").Value = 10

$ws1.Range("A60").Copy() | Out-Null
$ws2.Range("A40").PasteSpecial(-4122) | Out-Null
$ws2.Range("A40").Value = "FSB-3_L_Red Clover"
$ws2.Range("B40").Value = 2.5

$ws1.Range("A60").Copy() | Out-Null
$ws2.Range("A41").PasteSpecial(-4122) | Out-Null
$ws2.Range("A41").Value = "FSB-3_H_Oriental Mustard"
$ws2.Range("B41").Value = 20

$ws1.Range("A60").Copy() | Out-Null
$ws2.Range("A42").PasteSpecial(-4122) | Out-Null
$ws2.Range("A42").Value = "FSB-3_O_Control"
$ws2.Range("B42").Value = 0

$ws1.Range("A60").Copy() | Out-Null
$ws2.Range("A43").PasteSpecial(-4122) | Out-Null
$ws2.Range("A43").Value = "FSB-3_M_Alfalfa"
$ws2.Range("B43").Value = 5

$ws1.Range("A60").Copy() | Out-Null
$ws2.Range("A44").PasteSpecial(-4122) | Out-Null
$ws2.Range("A44").Value = "FSB-3_F_Mustard"
$ws2.Range("B44").Value = 30

$ws1.Range("A60").Copy() | Out-Null
$ws2.Range("A45").PasteSpecial(-4122) | Out-Null
$ws2.Range("A45").Value = "FSB-3_P_Control"
$ws2.Range("B45").Value = 0

$ws1.Range("A60").Copy() | Out-Null
$ws2.Range("A46").PasteSpecial(-4122) | Out-Null
$ws2.Range("A46").Value = "FSB-3_G_Brown Mustard"
$ws2.Range("B46").Value = 25

$ws1.Range("A60").Copy() | Out-Null
$ws2.Range("A47").PasteSpecial(-4122) | Out-Null
$ws2.Range("A47").Value = "FSB-3_I_Arugula"
$ws2.Range("B47").Value = 2.5

$ws1.Range("A60").Copy() | Out-Null
$ws2.Range("A48").PasteSpecial(-4122) | Out-Null
$ws2.Range("A48").Value = "FSB-3_K_Kale"
$ws2.Range("B48").Value = 5

$ws1.Range("A60").Copy() | Out-Null
$ws2.Range("A49").PasteSpecial(-4122) | Out-Null
$ws2.Range("A49").Value = "FSB-3_D_Collard"
$ws2.Range("B49").Value = 20

$ws1.Range("A89").Copy() | Out-Null
$ws2.Range("A50").PasteSpecial(-4122) | Out-Null
$ws2.Range("A50").Value = "FSB-4_N_Control"
$ws2.Range("B50").Value = 0

$ws1.Range("A89").Copy() | Out-Null
$ws2.Range("A51").PasteSpecial(-4122) | Out-Null
$ws2.Range("A51").Value = "FSB-4_M_Alfalfa"
$ws2.Range("B51").Value = 2.5

$ws1.Range("A89").Copy() | Out-Null
$ws2.Range("A52").PasteSpecial(-4122) | Out-Null
$ws2.Range("A52").Value = "FSB-4_P_Control"
$ws2.Range("B52").Value = 0

$ws1.Range("A89").Copy() | Out-Null
$ws2.Range("A53").PasteSpecial(-4122) | Out-Null
$ws2.Range("A53").Value = "FSB-4_B_Turnip"
$ws2.Range("B53").Value = 5

$ws1.Range("A89").Copy() | Out-Null
$ws2.Range("A54").PasteSpecial(-4122) | Out-Null
$ws2.Range("A54").Value = "FSB-4_E_Broadleaf Mustard"
$ws2.Range("B54").Value = 20

$ws1.Range("A89").Copy() | Out-Null
$ws2.Range("A55").PasteSpecial(-4122) | Out-Null
$ws2.Range("A55").Value = "FSB-4_F_Mustard"
$ws2.Range("B55").Value = 25

$ws1.Range("A89").Copy() | Out-Null
$ws2.Range("A56").PasteSpecial(-4122) | Out-Null
$ws2.Range("A56").Value = "FSB-4_H_Oriental Mustard"
$ws2.Range("B56").Value = 20

$ws1.Range("A89").Copy() | Out-Null
$ws2.Range("A57").PasteSpecial(-4122) | Out-Null
$ws2.Range("A57").Value = "FSB-4_C_Rapeseed"
$ws2.Range("B57").Value = 5

$ws1.Range("A89").Copy() | Out-Null
$ws2.Range("A58").PasteSpecial(-4122) | Out-Null
$ws2.Range("A58").Value = "FSB-4_D_Collard"
$ws2.Range("B58").Value = 30

$ws1.Range("A89").Copy() | Out-Null
$ws2.Range("A59").PasteSpecial(-4122) | Out-Null
$ws2.Range("A59").Value = "FSB-4_J_Winter Camelina"
$ws2.Range("B59").Value = 0

$ws1.Range("A89").Copy() | Out-Null
$ws2.Range("A60").PasteSpecial(-4122) | Out-Null
$ws2.Range("A60").Value = "FSB-4_O_Control"
$ws2.Range("B60").Value = 0

$ws1.Range("A89").Copy() | Out-Null
$ws2.Range("A61").PasteSpecial(-4122) | Out-Null
$ws2.Range("A61").Value = "FSB-4_L_Red Clover"
$ws2.Range("B61").Value = 5

$ws1.Range("A89").Copy() | Out-Null
$ws2.Range("A62").PasteSpecial(-4122) | Out-Null
$ws2.Range("A62").Value = "FSB-4_I_Arugula"
$ws2.Range("B62").Value = 2.5

$ws1.Range("A89").Copy() | Out-Null
$ws2.Range("A63").PasteSpecial(-4122) | Out-Null
$ws2.Range("A63").Value = "FSB-4_K_Kale"
$ws2.Range("B63").Value = 15

$ws1.Range("A89").Copy() | Out-Null
$ws2.Range("A64").PasteSpecial(-4122) | Out-Null
$ws2.Range("A64").Value = "FSB-4_G_Brown Mustard"
$ws2.Range("B64").Value = 7.5

$ws1.Range("A89").Copy() | Out-Null
$ws2.Range("A65").PasteSpecial(-4122) | Out-Null
$ws2.Range("A65").Value = "FSB-4_A_Radish"
$ws2.Range("B65").Value = 2.5

# Select the full populated range and make this sheet the active tab,
# matching the saved view state.
$ws2.Range("A1:B65").Select() | Out-Null
$ws2.Activate() | Out-Null